$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: add "Total" as the new last column (W1)
$ws.Range("W1").Value = "Total"

# Row totals for existing category rows 2-6 (column W)
$ws.Range("W2").Value = 1944
$ws.Range("W3").Value = 199
$ws.Range("W4").Value = 1074
$ws.Range("W5").Value = 280
$ws.Range("W6").Value = 1253

# New row 7: "Outros" category
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 157
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = 14
$ws.Range("E7").Value = 46
$ws.Range("F7").Value = 101
$ws.Range("G7").Value = 92
$ws.Range("H7").Value = 98
$ws.Range("I7").Value = 101
$ws.Range("J7").Value = 106
$ws.Range("K7").Value = 124
$ws.Range("L7").Value = 101
$ws.Range("M7").Value = 108
$ws.Range("N7").Value = 113
$ws.Range("O7").Value = 114
$ws.Range("P7").Value = 115
$ws.Range("Q7").Value = 130
$ws.Range("R7").Value = 151
$ws.Range("S7").Value = 144
$ws.Range("T7").Value = 89
$ws.Range("U7").Value = 35
$ws.Range("V7").Value = 10
$ws.Range("W7").Value = 1960

# New row 8: "Total" row (column totals across all category rows)
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 174
$ws.Range("C8").Value = 16
$ws.Range("D8").Value = 22
$ws.Range("E8").Value = 55
$ws.Range("F8").Value = 116
$ws.Range("G8").Value = 112
$ws.Range("H8").Value = 129
$ws.Range("I8").Value = 163
$ws.Range("J8").Value = 224
$ws.Range("K8").Value = 292
$ws.Range("L8").Value = 358
$ws.Range("M8").Value = 470
$ws.Range("N8").Value = 522
$ws.Range("O8").Value = 522
$ws.Range("P8").Value = 691
$ws.Range("Q8").Value = 755
$ws.Range("R8").Value = 820
$ws.Range("S8").Value = 679
$ws.Range("T8").Value = 394
$ws.Range("U8").Value = 157
$ws.Range("V8").Value = 39
$ws.Range("W8").Value = 6710
